$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing rows (old rows 8-10, the "MuSCs-sending" TPM block)
$ws.Range("A8:T10").EntireRow.Delete()

# New TPM values for rows 2-7 (sending cluster FAPs/MuSCs x target cluster ECs/FAPs/MuSCs)
# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.071953
$ws.Range("H2").Value = 126.215859
$ws.Range("I2").Value = 0.978774012990499
$ws.Range("J2").Value = 0.978774012990499
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 5186.657390057811
$ws.Range("R2").Value = 46679.91651052029
$ws.Range("S2").Value = 0.6109090492400876
$ws.Range("T2").Value = 0.6109090492400875

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.071953
$ws.Range("H3").Value = 126.215859
$ws.Range("I3").Value = 0.978774012990499
$ws.Range("J3").Value = 0.978774012990499
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 1988.322952977258
$ws.Range("R3").Value = 17894.90657679532
$ws.Range("S3").Value = 0.2341940856000981
$ws.Range("T3").Value = 0.2341940856000981

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.071953
$ws.Range("H4").Value = 126.215859
$ws.Range("I4").Value = 0.978774012990499
$ws.Range("J4").Value = 0.978774012990499
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 1134.874411921455
$ws.Range("R4").Value = 10213.8697072931
$ws.Range("S4").Value = 0.1336708781503134
$ws.Range("T4").Value = 0.1336708781503133

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.912385
$ws.Range("H5").Value = 2.737155
$ws.Range("I5").Value = 0.021225987009501
$ws.Range("J5").Value = 0.021225987009501
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 112.4794088552983
$ws.Range("R5").Value = 1012.314679697685
$ws.Range("S5").Value = 0.01324835699666515
$ws.Range("T5").Value = 0.01324835699666515

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.912385
$ws.Range("H6").Value = 2.737155
$ws.Range("I6").Value = 0.021225987009501
$ws.Range("J6").Value = 0.021225987009501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("Q6").Value = 43.11936832245833
$ws.Range("R6").Value = 388.074314902125
$ws.Range("S6").Value = 0.005078803230034163
$ws.Range("T6").Value = 0.005078803230034163

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.912385
$ws.Range("H7").Value = 2.737155
$ws.Range("I7").Value = 0.021225987009501
$ws.Range("J7").Value = 0.021225987009501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("Q7").Value = 24.61122711182334
$ws.Range("R7").Value = 221.50104400641
$ws.Range("S7").Value = 0.002898826782801684
$ws.Range("T7").Value = 0.002898826782801683
